# DEAN import format: add the missing "Course" record for PSY180
# (Course_Section already references course_section_code "PSY180" /
# course_id 1395954, but the Course sheet had no matching row), and
# leave the "Course" sheet as the active tab/selection, matching the
# author's saved workbook state.

$wb = $excel.ActiveWorkbook

# Record the selection left behind on Course_Section (it stops being the
# active sheet, but Excel still remembers where the cursor was on it).
$courseSection = $wb.Worksheets.Item("Course_Section")
$courseSection.Activate()
$null = $courseSection.Range("G22").Select()

# Course sheet: append the course row referenced by Course_Section.
$course = $wb.Worksheets.Item("Course")
$course.Range("A2").Value = 1395954
$course.Range("B2").Value = "PSY180"
$course.Range("C2").Value = "Interpersonal Effectiveness"

# Course becomes the active sheet/cell in the saved workbook.
$course.Activate()
$null = $course.Range("A4").Select()
